$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 72.2
$ws.Range("I39").Value = 14.769231
$ws.Range("J39").Value = 178.85715
$ws.Range("K39").Value = 44.307693
$ws.Range("L39").Value = 536.5714499999999
$ws.Range("M39").Value = 251.692307
$ws.Range("N39").Value = -1128.57145

$ws.Range("H40").Value = 4616.222
$ws.Range("I40").Value = 2986.75
$ws.Range("J40").Value = 5919.8
$ws.Range("K40").Value = 2986.75
$ws.Range("L40").Value = 5919.8
$ws.Range("M40").Value = -2811.75
$ws.Range("N40").Value = -6269.8

$ws.Range("H62").Value = 5445.4546
$ws.Range("I62").Value = 6402.25
$ws.Range("K62").Value = 6402.25
$ws.Range("M62").Value = -5778.25

$ws.Range("H65").Value = 5445.4546
$ws.Range("I65").Value = 6402.25
$ws.Range("K65").Value = 32011.25
$ws.Range("M65").Value = -28891.25

$ws.Range("H100").Value = 2541
$ws.Range("I100").Value = 2541
$ws.Range("K100").Value = 2541
$ws.Range("M100").Value = -2000

$ws.Range("H132").Value = 5750.7334
$ws.Range("I132").Value = 2374.25
$ws.Range("K132").Value = 7122.75
$ws.Range("M132").Value = -4592.75

$ws.Range("H137").Value = 1795.3636
$ws.Range("I137").Value = 1392.7142
$ws.Range("K137").Value = 4178.142599999999
$ws.Range("M137").Value = -1628.142599999999

$ws.Range("H138").Value = 5872.7896
$ws.Range("I138").Value = 5039
$ws.Range("K138").Value = 15117
$ws.Range("M138").Value = -9977

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 13499.5
$ws.Range("J61").Value = 7000
$ws.Range("L61").Value = 7000
$ws.Range("N61").Value = -7424

$ws.Range("H74").Value = 1472.1428
$ws.Range("I74").Value = 1416.8182
$ws.Range("K74").Value = 1416.8182
$ws.Range("M74").Value = -542.8181999999999

$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20676

$ws.Range("H77").Value = 1472.1428
$ws.Range("I77").Value = 1416.8182
$ws.Range("K77").Value = 7084.090999999999
$ws.Range("M77").Value = -2716.090999999999

$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22340

$ws.Range("H88").Value = 4518.364
$ws.Range("I88").Value = 3502
$ws.Range("J88").Value = 5365.3335
$ws.Range("K88").Value = 3502
$ws.Range("L88").Value = 5365.3335
$ws.Range("M88").Value = -3096
$ws.Range("N88").Value = -6177.3335

$ws.Range("H91").Value = 4518.364
$ws.Range("I91").Value = 3502
$ws.Range("J91").Value = 5365.3335
$ws.Range("K91").Value = 3502
$ws.Range("L91").Value = 5365.3335
$ws.Range("M91").Value = -2098
$ws.Range("N91").Value = -8173.3335

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H132").Value = 3570.1428
$ws.Range("I132").Value = 3165.1667
$ws.Range("J132").Value = 6000
$ws.Range("K132").Value = 9495.500100000001
$ws.Range("L132").Value = 18000
$ws.Range("M132").Value = -6965.500100000001
$ws.Range("N132").Value = -23060

$ws.Range("H136").Value = 13499.5
$ws.Range("J136").Value = 7000
$ws.Range("L136").Value = 21000
$ws.Range("N136").Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H60").Value = 39259.332
$ws.Range("J60").Value = 39259.332
$ws.Range("L60").Value = 39259.332
$ws.Range("N60").Value = -40457.332

$ws.Range("H134").Value = 2990.2
$ws.Range("I134").Value = 1652.3334
$ws.Range("K134").Value = 4957.0002
$ws.Range("M134").Value = -2422.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2905.889
$ws.Range("I31").Value = 2905.889
$ws.Range("K31").Value = 2905.889
$ws.Range("M31").Value = -2610.889

$ws.Range("H34").Value = 2905.889
$ws.Range("I34").Value = 2905.889
$ws.Range("K34").Value = 2905.889
$ws.Range("M34").Value = -2703.889

$ws.Range("H105").Value = 2077
$ws.Range("I105").Value = 1398
$ws.Range("K105").Value = 1398
$ws.Range("M105").Value = 349

$ws.Range("H122").Value = 920.6667
$ws.Range("I122").Value = 920.6667
$ws.Range("K122").Value = 2762.0001
$ws.Range("M122").Value = -312.0001000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1387.5
$ws.Range("I131").Value = 750
$ws.Range("K131").Value = 2250
$ws.Range("M131").Value = 2790

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 42987.5
$ws.Range("I62").Value = 42987.5
$ws.Range("K62").Value = 42987.5
$ws.Range("M62").Value = -42301.5

$ws.Range("H65").Value = 42987.5
$ws.Range("I65").Value = 42987.5
$ws.Range("K65").Value = 128962.5
$ws.Range("M65").Value = -125530.5

$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 8000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 24000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -21530
$ws.Range("N126").ClearContents()

$ws.Range("H138").Value = 65000
$ws.Range("J138").Value = 65000
$ws.Range("L138").Value = 65000
$ws.Range("N138").Value = -75280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1074.25
$ws.Range("I16").Value = 1269.6666
$ws.Range("K16").Value = 1269.6666
$ws.Range("M16").Value = -1099.6666

$ws.Range("H46").Value = 3172.7273
$ws.Range("I46").Value = 2950
$ws.Range("J46").Value = 3222.2222
$ws.Range("K46").Value = 2950
$ws.Range("L46").Value = 3222.2222
$ws.Range("M46").Value = -2762
$ws.Range("N46").Value = -3598.2222

$ws.Range("H97").Value = 9475
$ws.Range("J97").Value = 9475
$ws.Range("L97").Value = 9475
$ws.Range("N97").Value = -11457

$ws.Range("H100").Value = 1996.3334
$ws.Range("I100").Value = 1595.6
$ws.Range("K100").Value = 1595.6
$ws.Range("M100").Value = -1054.6

$ws.Range("H122").Value = 6245.8335
$ws.Range("I122").Value = 4342.857
$ws.Range("J122").Value = 8910
$ws.Range("K122").Value = 13028.571
$ws.Range("L122").Value = 26730
$ws.Range("M122").Value = -10578.571
$ws.Range("N122").Value = -31630

$ws.Range("H132").Value = 5102.4287
$ws.Range("I132").Value = 4619.6665
$ws.Range("K132").Value = 13858.9995
$ws.Range("M132").Value = -11328.9995

$ws.Range("H136").Value = 4519.5
$ws.Range("I136").Value = 3719.3333
$ws.Range("K136").Value = 11157.9999
$ws.Range("M136").Value = -8607.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 29995
$ws.Range("I52").Value = 29995
$ws.Range("K52").Value = 29995
$ws.Range("M52").Value = -29769

$ws.Range("H132").Value = 1982.375
$ws.Range("I132").Value = 1982.375
$ws.Range("K132").Value = 5947.125
$ws.Range("M132").Value = -3417.125
